$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# 1) "Who are we, ... mud, godless ..." -> "Who are we- ... mud- godless ..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Who are we, creatures from the mud, godless now and doomed to roam, eternally searching for that which we will never find.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Who are we- creatures from the mud- godless now and doomed to roam, eternally searching for that which we will never find.",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# Helper: locate paragraphs by their exact text content.
# ---------------------------------------------------------------------------
function Find-ParagraphIndex($doc, [string]$text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7) -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 2) Move the _GoBack bookmark from the end of the "nature of things"
#    paragraph to the end of the "A Note" heading paragraph (just after the
#    run, before the paragraph mark).
# ---------------------------------------------------------------------------
$noteIdx = Find-ParagraphIndex $d "A Note"
$pNote = $d.Paragraphs($noteIdx)
$xmlNote = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="3"/><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-GB"/></w:rPr><w:t>A Note</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$pNote.Range.InsertXML($xmlNote) | Out-Null

$natureBodyText = "It is the nature of things that every event should be balanced by another. Where there is force, there is friction. Where there is unity, there is discord. Where there is love there is hate. Where there is light there is darkness. And where there is life, there is also death. It is the nature of things, and we are mere subjects to the will of nature."
$natureIdx = Find-ParagraphIndex $d $natureBodyText
$pNature = $d.Paragraphs($natureIdx)
$xmlNature = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-GB"/></w:rPr><w:t>' + $natureBodyText + '</w:t></w:r></w:p>'
$pNature.Range.InsertXML($xmlNature) | Out-Null

# ---------------------------------------------------------------------------
# 3) Append the new "Apocrypha" section after the "nature of things" body
#    paragraph: blank line, heading, blank line, body quote paragraph.
# ---------------------------------------------------------------------------
$xmlBlank = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p>'

$natureIdx = Find-ParagraphIndex $d $natureBodyText
$pNature = $d.Paragraphs($natureIdx)
$r = $pNature.Range
$r.Collapse(0) | Out-Null
$r.InsertParagraphAfter() | Out-Null
$blank1 = $d.Paragraphs($natureIdx + 1)
$blank1.Range.InsertXML($xmlBlank) | Out-Null

$blank1.Range.Collapse(0) | Out-Null
$blank1.Range.InsertParagraphAfter() | Out-Null
$headPara = $d.Paragraphs($natureIdx + 2)
$xmlHead = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="3"/><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-GB"/></w:rPr><w:t>Apocrypha</w:t></w:r></w:p>'
$headPara.Range.InsertXML($xmlHead) | Out-Null

$headPara.Range.Collapse(0) | Out-Null
$headPara.Range.InsertParagraphAfter() | Out-Null
$blank2 = $d.Paragraphs($natureIdx + 3)
$blank2.Range.InsertXML($xmlBlank) | Out-Null

$blank2.Range.Collapse(0) | Out-Null
$blank2.Range.InsertParagraphAfter() | Out-Null
$quotePara = $d.Paragraphs($natureIdx + 4)
$quoteText = @"
“…and as I stood on the shore I saw a beast rise up from the water and as it fell from it’s scaly back it said to me ‘I am the one who rises up at the end of days and those who see me shall know that I am death’. And with that the beast rose up and walked upon the land, and I knew that it was true, that the end of days were upon us, that the one of whom the Necromancer spake would arrive and would signal the death knell of the world.” - Unknown Doomsayer
"@
$xmlQuote = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-GB"/></w:rPr><w:t>' + $quoteText + '</w:t></w:r></w:p>'
$quotePara.Range.InsertXML($xmlQuote) | Out-Null

Write-Output "All edits applied."
